# Updated cryptos list on Wed May  1 10:40:19 UTC 2024 with GitHub Actions
#
# This script refreshes the "Price" (column D) and "Volume(1h)" (column E)
# figures for the crypto table on the active sheet, and also fixes the
# Filecoin/Mantle rows (35/36), whose coin name, link, price and volume
# all changed.
#
# A handful of the new "Price" strings look like plain numbers (e.g.
# "550.15", "4.72", "2.00") even though the column stores them as text
# (others, like "57.449.20", have two '.' separators and are never
# auto-converted). To keep those cells as text instead of letting Excel
# silently coerce them to numeric values, we force NumberFormat="@" on
# just those cells before writing their new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Force "text" format on Price cells whose new value would otherwise
#     be auto-converted to a number ---------------------------------------
$ws.Range("D5").NumberFormat  = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# --- Row 2 : Bitcoin ------------------------------------------------------
$ws.Range('D2').Value = '57.449.20'
$ws.Range('E2').Value = '  -6.83%  '

# --- Row 3 : Ethereum -----------------------------------------------------
$ws.Range('D3').Value = '2.885.97'
$ws.Range('E3').Value = '  -5.14%  '

# --- Row 4 : TetherUSD -----------------------------------------------------
$ws.Range('E4').Value = '  +0.02%  '

# --- Row 5 : BNB -----------------------------------------------------------
$ws.Range('D5').Value = '550.15'
$ws.Range('E5').Value = '  -5.27%  '

# --- Row 6 : Solana ---------------------------------------------------------
$ws.Range('E6').Value = '  -6.93%  '

# --- Row 7 : USDC ------------------------------------------------------------
$ws.Range('E7').Value = '  +0.03%  '

# --- Row 8 : LidoStakedEther --------------------------------------------------
$ws.Range('D8').Value = '2.881.50'
$ws.Range('E8').Value = '  -5.26%  '

# --- Row 9 : XRP ---------------------------------------------------------------
$ws.Range('E9').Value = '  -1.62%  '

# --- Row 10 : Dogecoin -----------------------------------------------------------
$ws.Range('E10').Value = '  -9.81%  '

# --- Row 11 : Toncoin -------------------------------------------------------------
$ws.Range('D11').Value = '4.72'
$ws.Range('E11').Value = '  -9.44%  '

# --- Row 12 : Cardano --------------------------------------------------------------
$ws.Range('D12').Value = '0.431'
$ws.Range('E12').Value = '  -1.52%  '

# --- Row 13 : ShibaInu ---------------------------------------------------------------
$ws.Range('E13').Value = '  -9.41%  '

# --- Row 14 : Avalanche ----------------------------------------------------------------
$ws.Range('D14').Value = '31.39'
$ws.Range('E14').Value = '  -5.88%  '

# --- Row 15 : TRON -----------------------------------------------------------------------
$ws.Range('E15').Value = '  -1.06%  '

# --- Row 16 : WrappedliquidstakedEther2.0 -------------------------------------------------
$ws.Range('D16').Value = '3.359.05'
$ws.Range('E16').Value = '  -5.18%  '

# --- Row 17 : WrappedEther -----------------------------------------------------------------
$ws.Range('D17').Value = '2.881.02'
$ws.Range('E17').Value = '  -5.53%  '

# --- Row 18 : Polkadot ----------------------------------------------------------------------
$ws.Range('D18').Value = '6.47'
$ws.Range('E18').Value = '  +1.83%  '

# --- Row 19 : WrappedBTC --------------------------------------------------------------------
$ws.Range('D19').Value = '57.394.40'
$ws.Range('E19').Value = '  -6.90%  '

# --- Row 20 : BitcoinCash -------------------------------------------------------------------
$ws.Range('D20').Value = '406.97'
$ws.Range('E20').Value = '  -8.78%  '

# --- Row 21 : Chainlink ---------------------------------------------------------------------
$ws.Range('D21').Value = '12.75'
$ws.Range('E21').Value = '  -5.24%  '

# --- Row 22 : Polygon -----------------------------------------------------------------------
$ws.Range('D22').Value = '0.650'
$ws.Range('E22').Value = '  -2.76%  '

# --- Row 23 : Uniswap -----------------------------------------------------------------------
$ws.Range('E23').Value = '  -8.33%  '

# --- Row 24 : InternetComputer(DFINITY) ------------------------------------------------------
$ws.Range('D24').Value = '12.54'
$ws.Range('E24').Value = '  -2.07%  '

# --- Row 25 : Litecoin ------------------------------------------------------------------------
$ws.Range('D25').Value = '76.58'
$ws.Range('E25').Value = '  -5.16%  '

# --- Row 27 : FirstDigitalUSD -------------------------------------------------------------------
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.08%  '

# --- Row 28 : PancakeSwap -----------------------------------------------------------------------
$ws.Range('E28').Value = '  -4.06%  '

# --- Row 29 : ImmutableX ------------------------------------------------------------------------
$ws.Range('E29').Value = '  -5.10%  '

# --- Row 30 : RenderToken -----------------------------------------------------------------------
$ws.Range('D30').Value = '7.07'
$ws.Range('E30').Value = '  -4.57%  '

# --- Row 31 : NEARProtocol ----------------------------------------------------------------------
$ws.Range('D31').Value = '6.00'
$ws.Range('E31').Value = '  -7.12%  '

# --- Row 32 : EthereumClassic -------------------------------------------------------------------
$ws.Range('D32').Value = '24.48'
$ws.Range('E32').Value = '  -5.25%  '

# --- Row 33 : Hedera ----------------------------------------------------------------------------
$ws.Range('D33').Value = '0.0938'
$ws.Range('E33').Value = '  -3.01%  '

# --- Row 34 : Stacks ----------------------------------------------------------------------------
$ws.Range('D34').Value = '2.00'
$ws.Range('E34').Value = '  -13.54%  '

# --- Row 35 : was Filecoin, is now Mantle --------------------------------------------------------
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').Value = '0.893'
$ws.Range('E35').Value = '  -7.78%  '

# --- Row 36 : was Mantle, is now Filecoin --------------------------------------------------------
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '5.32'
$ws.Range('E36').Value = '  -6.47%  '

# --- Row 37 : OKB -----------------------------------------------------------------------------------
$ws.Range('D37').Value = '48.25'
$ws.Range('E37').Value = '  -4.02%  '

# --- Row 38 : Cosmos --------------------------------------------------------------------------------
$ws.Range('D38').Value = '8.31'
$ws.Range('E38').Value = '  +5.24%  '

# --- Row 39 : PEPE ----------------------------------------------------------------------------------
$ws.Range('E39').Value = '  -12.34%  '

# --- Row 40 : VeChain -------------------------------------------------------------------------------
$ws.Range('D40').Value = '0.0342'
$ws.Range('E40').Value = '  -8.31%  '

# --- Row 41 : Kaspa ----------------------------------------------------------------------------------
$ws.Range('D41').Value = '0.104'
$ws.Range('E41').Value = '  -4.22%  '

# --- Row 42 : Bittensor -------------------------------------------------------------------------------
$ws.Range('D42').Value = '361.00'
$ws.Range('E42').Value = '  -4.29%  '

# --- Row 43 : Maker -----------------------------------------------------------------------------------
$ws.Range('D43').Value = '2.601.67'
$ws.Range('E43').Value = '  -3.12%  '

# --- Row 44 : USDe -------------------------------------------------------------------------------------
$ws.Range('E44').Value = '  -0.05%  '

# --- Row 45 : dogwifhat --------------------------------------------------------------------------------
$ws.Range('E45').Value = '  -8.20%  '

# --- Row 46 : Monero -----------------------------------------------------------------------------------
$ws.Range('D46').Value = '117.76'
$ws.Range('E46').Value = '  -4.19%  '

# --- Row 47 : TheGraph ----------------------------------------------------------------------------------
$ws.Range('E47').Value = '  -4.80%  '

# --- Row 48 : Stellar -----------------------------------------------------------------------------------
$ws.Range('E48').Value = '  -1.97%  '

# --- Row 49 : Fetch.AI ----------------------------------------------------------------------------------
$ws.Range('E49').Value = '  -4.37%  '

# --- Row 50 : InjectiveProtocol -------------------------------------------------------------------------
$ws.Range('D50').Value = '22.08'
$ws.Range('E50').Value = '  -7.12%  '

# --- Row 51 : ThetaToken --------------------------------------------------------------------------------
$ws.Range('E51').Value = '  -5.66%  '
